$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(514, '피부색이 창백하다', 'pale skin', '205'),
    @(515, '피부가 희다', 'fair skin', '205'),
    @(516, '피부가 햇볕에 탔다', 'tan skin', '205'),
    @(517, '피부색이 검다', 'dark skin', '205'),
    @(518, '얼굴이 둥글다', 'a round face', '205'),
    @(519, '얼굴이 달걀형이다', 'an oval-shaped face', '205'),
    @(520, '얼굴이 갸름하다', 'a thin face', '205'),
    @(521, '얼굴이 갸름하다', 'an oblong face', '205'),
    @(522, '얼굴이 사각형이다', 'a square face', '205'),
    @(523, '여드름, 뾰루지', 'pimple', '205'),
    @(524, '다크서클', 'dark circles', '205'),
    @(525, '주근깨, 기미', 'freckle', '205'),
    @(526, '주름', 'wrinkle', '205'),
    @(527, '여드름, 뾰루지', 'acne', '205'),
    @(528, '건성 피부', 'dry skin', '205'),
    @(529, '점', 'mole', '205'),
    @(530, '보조개', 'dimple', '205'),
    @(531, '지성 피부', 'oily skin', '205'),
    @(532, '머리가 길다', 'long hair', '205'),
    @(533, '직모다', 'straight hair', '205'),
    @(534, '머리가 짧다', 'short hair', '205'),
    @(535, '머리가 어깨까지 오다', 'shoulder-length hair', '205'),
    @(536, '웨이브가 있다', 'wavy hair', '205'),
    @(537, '곱슬머리다', 'curly hair', '205'),
    @(538, '흰머리가 있다', 'gray hair', '205'),
    @(539, '머리를 뒤로 묶다', 'wear a pony tail', '205'),
    @(540, '머리를 자르다', 'have one''s hair cut', '205'),
    @(541, '파마하다', 'have one''s hair permed', '205'),
    @(542, '염색하다', 'have one''s hair dyed', '205'),
    @(543, '머리가 빠지다', 'lose one''s hair', '205'),
    @(544, '대머리다', 'be bald', '205'),
    @(545, '턱수염이 있다', 'have beard', '205'),
    @(546, '콧수염이 있다', 'have mustache', '205'),
    @(547, '구레나룻이 있다', 'have sideburns', '205'),
    @(548, '숨을 쉬다', 'breathe', '206'),
    @(549, '숨을 참다', 'hold one''s breath', '206'),
    @(550, '한숨을 쉬다, 한숨', 'sigh', '206'),
    @(551, '하품하다', 'yawn', '206'),
    @(552, '기침하다', 'cough', '206'),
    @(553, '재채기하다', 'sneeze', '206'),
    @(554, '딸꾹질하다', 'hiccup', '206'),
    @(555, '딸꾹질', 'hiccups', '206'),
    @(556, '윙크하다 (한쪽)', 'wink', '206'),
    @(557, '눈을 깜박이다 (양쪽)', 'blink', '206'),
    @(558, '미소 짓다', 'smile', '206'),
    @(559, '소리 내어 웃다', 'laugh', '206'),
    @(560, '찡그리다', 'frown', '206'),
    @(561, '울다', 'cry', '206'),
    @(562, '눈물을 흘리다', 'weep', '206'),
    @(563, '뒤를 돌아보다', 'look back', '206'),
    @(564, '고개를 끄덕이다', 'nod', '206'),
    @(565, '고개를 젓다', 'snake one''s head', '206'),
    @(566, '고개를 숙이다', 'lower one''s head', '206'),
    @(567, '콧물이 나다', 'one''s nose runs', '206'),
    @(568, '코를 풀다', 'blow one''s nose', '206'),
    @(569, '코를 닦다', 'wipe one''s nose', '206'),
    @(570, '코딱지를 파다', 'pick one''s nose', '206'),
    @(571, '침을 뱉다', 'spit', '206'),
    @(572, '가래를 뱉다', 'spit out phlegm', '206'),
    @(573, '속삭이다', 'whisper', '206'),
    @(574, '소리치다, 고함치다', 'shout', '206'),
    @(575, '손을 들다', 'raise one''s hand', '207'),
    @(576, '손을 흔들다', 'wave', '207'),
    @(577, '박수를 치다', 'clap one''s hand', '207'),
    @(578, '~와 악수하다', 'shake hands with ~', '207'),
    @(579, '팔짱을 끼다', 'fold one''s arms', '207'),
    @(580, '들다, 휴대하다, 나르다', 'carry', '207'),
    @(581, '집다, 집어서 들어 올리다', 'pick up', '207'),
    @(582, '만지다', 'touch', '207'),
    @(583, '가리키다', 'point', '207'),
    @(584, '가리키다', 'point at', '207'),
    @(585, '쥐고 있다, 잡고 있다', 'hold', '207'),
    @(586, '치다, 때리다', 'hit', '207'),
    @(587, '위로 들어 올리다', 'lift', '207'),
    @(588, '던지다', 'throw', '207'),
    @(589, '(움직이는 물체를) 잡다, 받다', 'catch', '207'),
    @(590, '당기다', 'pull', '207'),
    @(591, '밀다', 'push', '207'),
    @(592, '손으로 짜다', 'squeeze', '207'),
    @(593, '비틀다, 구부리다', 'twist', '207'),
    @(594, '눕다', 'lie', '207'),
    @(595, '엎드리다', 'lie on one''s face', '207'),
    @(596, '엎드리다', 'lie on one''s stomach', '207'),
    @(597, '일어서다', 'stand up', '207'),
    @(598, '쓰러지다, 넘어지다', 'fall down', '207'),
    @(599, '허리굽혀 인사하다', 'bow', '207'),
    @(600, '고개숙여 인사하다', 'bow', '207'),
    @(601, '어깨를 으쓱하다', 'shrug', '207'),
    @(602, '몸을 떨다', 'shiver', '207'),
    @(603, '안다, 포옹하다', 'hug', '207'),
    @(604, '안다, 포옹하다', 'embrace', '207'),
    @(605, '걷다', 'walk', '207'),
    @(606, '뛰다, 달리다', 'run', '207'),
    @(607, '점프하다', 'jump', '207'),
    @(608, '무릎을 꿇다', 'kneel', '207'),
    @(609, '무릎을 꿇다', 'kneel down', '207'),
    @(610, '발로 차다', 'kick', '207'),
    @(611, '(엎드려) 기다', 'crawl', '207'),
    @(612, '올라가다, 오르다', 'climb', '207'),
    @(613, '무릎을 구부리다', 'bend one''s knees', '207'),
    @(614, '발끝으로 살금살금 걷다', 'tiptoe', '207'),
    @(615, '발끝으로 살금살금 걷다', 'walk on tiptoe', '207')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

[void]$ws.Range("D615").Select()
try { $excel.ActiveWindow.ScrollRow = 599 } catch {}
try { $ws.PageSetup.Orientation = 1 } catch {}

Write-Host "Added $($data.Count) rows"
